$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: new data for IP 85.104.3.240 ---
# The order in which brand-new text values are first written controls the
# order they are appended to the shared string table, so the assignments
# below are deliberately sequenced to match the source data.

# D4: hyperlink cell (add the relationship first, then restore the same
# cell formatting used by the other hyperlink cells D2/D3 so no extra
# style gets introduced).
$ws.Hyperlinks.Add($ws.Range("D4"), "https://www.virustotal.com/gui/ip-address/85.104.3.240/detection") | Out-Null
$ws.Range("D4").Value = "https://www.virustotal.com/gui/ip-address/85.104.3.240/detection"
$ws.Range("D4").Style = "Hyperlink"

# E4: last_analysis_stats
$ws.Range("E4").Value = "{'harmless': 56, 'malicious': 12, 'suspicious': 1, 'undetected': 21, 'timeout': 0}"

# F4: Country
$ws.Range("F4").Value = "Turkey"

# C4: Status (reuses the existing "Malicious" shared string)
$ws.Range("C4").Value = "Malicious"

# G4/H4: dates, using the same number format as the rows above
$ws.Range("G2").Copy($ws.Range("G4"))
$ws.Range("G4").Value = 45306.97109953704
$ws.Range("H2").Copy($ws.Range("H4"))
$ws.Range("H4").Value = 45340.0480787037

# I4: AS_Owner
$ws.Range("I4").Value = "Turk Telekom"
